# Update "Documentação de Convênios Saída" workbook:
#  - sheet 1 (Documentaçãometadados): update the "Descrição" value (B4)
#  - sheet 2 (Dicionário): fill in the data-dictionary rows (A2:D22) describing
#    the columns of the "valores_indicados_emendas_impos" report, up to column U's
#    worth of report fields (cnpj gets a numeric "0" format).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet 1 - Documentaçãometadados: "Descrição" field gets an updated value
# ---------------------------------------------------------------------------
$ws1.Range("B4").Value = "Relatórios contendo o status e a execução das emendas parlamentares - Emendas 2020. "

# ---------------------------------------------------------------------------
# Sheet 2 - Dicionário: data dictionary rows for the report
# ---------------------------------------------------------------------------
$reportName = "valores_indicados_emendas_impos"

$rows = @(
    @{ Row=2;  B="responsavel";            C="Responsável pela emenda";       D="pode ser deputado etc" },
    @{ Row=3;  B="uo";                     C="Unidade Organizacional";        D=$null },
    @{ Row=4;  B="orgao";                  C="Órgão";                         D=$null },
    @{ Row=5;  B="n_sigcon";               C="Número de Referência SIGCON";   D=$null },
    @{ Row=6;  B="tipo_indicacao";         C="Tipo de Indicação";             D=$null },
    @{ Row=7;  B="ano_inciso";             C="Ano do Inciso";                 D=$null },
    @{ Row=8;  B="impositividade";         C="Impositividade";                D=$null },
    @{ Row=9;  B="municipio";              C="Município";                     D=$null },
    @{ Row=10; B="razao_social";           C="Razão Social";                  D="Razão social do beneficiário" },
    @{ Row=11; B="cnpj";                   C="CNPJ";                          D="CNPJ do beneficiário" },
    @{ Row=12; B="codesc";                 C="Código Escola";                 D="Código das escolas, intitutos e centros educacionais" },
    @{ Row=13; B="acao";                   C="Ação";                          D="Ação da LOA 2020" },
    @{ Row=14; B="grupo_despesa";          C="Grupo de Despesa";              D="Grupo de despesa LOA 2020" },
    @{ Row=15; B="genero";                 C="Gênero";                        D="Gênero da despesa" },
    @{ Row=16; B="categoria";              C="Categoria";                     D=$null },
    @{ Row=17; B="especificacao";          C="Especificação";                 D=$null },
    @{ Row=18; B="descricao";              C="Descrição";                     D=$null },
    @{ Row=19; B="tipo_aplicacao";         C="Tipo de Aplicação";             D=$null },
    @{ Row=20; B="vl_indicado";            C="Valor Indicado";                D=$null },
    @{ Row=21; B="data_indicacao_sigcon";  C="Data da Indicação no SIGCON";   D=$null },
    @{ Row=22; B="vl_informado";           C="Valor Informado";               D=$null }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws2.Cells.Item($rowNum, 1).Value = $reportName
    $ws2.Cells.Item($rowNum, 2).Value = $r.B
    $ws2.Cells.Item($rowNum, 3).Value = $r.C
    if ($r.D) {
        $ws2.Cells.Item($rowNum, 4).Value = $r.D
    }
}

# "cnpj" is a numeric identifier column -> integer number format (row 11)
$ws2.Range("B11").NumberFormat = "0"

# ---------------------------------------------------------------------------
# Selections, matching where the author last left the cursor in each sheet
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B4").Select()

$ws2.Activate()
$ws2.Range("A2:A22").Select()
